$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (copy formatting from the existing MOP_DESC header cell)
$ws.Range("F1").Value = "MOP_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New column data (MOP_DEF values, plain text like columns C/E)
$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"
$ws.Range("F3").Value = "['B is a process profile if there is some process c such that b is process profile of c. B is process profile of c holds when b is a proper occurrent part of c and there is some proper occurrent part d of c which has no parts in common with b and is mutually dependent on b and is such that b, c and d occupy the same temporal region. [BFO]']"
$ws.Range("F4").Value = "['B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]']"
$ws.Range("F5").Value = "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']"
$ws.Range("F6").Value = "['A subatomic particle is a material that is below the scale of an atom. [Allotrope]']"
$ws.Range("F7").Value = "[]"

$excel.CutCopyMode = $false
